$d = $word.ActiveDocument

# --- Edit 1: insert a new "Meta description" paragraph right after the title heading ---
$titlePara = $d.Paragraphs.Item(1)
$null = $titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience Dragon''s Fire, a stunning online slot game with 2D and 3D animations. Play for free on desktops and mobile devices.</w:t></w:r></w:p>'
$null = $metaPara.Range.InsertXML($metaXml)

# --- Edit 2: remove the duplicated bold title paragraph near the end, and replace the
#     remaining italic paragraph's text with the new image prompt ---
$count = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($count - 1)
$boldTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Text = "Prompt: Create a feature image for Dragon's Fire that showcases a happy Maya warrior with glasses amidst the dragon-themed slot machine. The image should be in a cartoon style and must be eye-catching to suit the game's mesmerizing graphics."
